$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.235.00'
$ws.Range("E2").Value = '  +1.15%  '

$ws.Range("D3").Value = '1.655.17'
$ws.Range("E3").Value = '  +1.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.02'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.00%  '

$ws.Range("E7").Value = '  +1.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0640'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.255'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0802'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.01%  '

$ws.Range("D12").Value = '1.693.30'
$ws.Range("E12").Value = '  +3.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.28'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.544'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.69'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.15%  '

$ws.Range("D16").Value = '0.0₃0763'
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").Value = '26.130.16'
$ws.Range("E17").Value = '  +0.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '194.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.17%  '

$ws.Range("E23").Value = '  +1.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.68%  '

$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0490'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.906'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.42%  '

$ws.Range("D36").Value = '1.140.95'
$ws.Range("E36").Value = '  +0.17%  '

$ws.Range("E37").Value = '  +0.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.535'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0157'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.803'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.68%  '

$ws.Range("D43").Value = '0.0₆0114'
$ws.Range("E43").Value = '  -1.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '56.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0524'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.420'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.12%  '

$ws.Range("E49").Value = '  +1.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0945'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.88%  '
